$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections -------------------------------------------------
# Employee name casing fix
$ws.Range("B5").Value = "John Vincent ramiso ramada"

# Month label fix
$ws.Range("B7").Value = "For the month of March - March"

# --- Row 12 time values collapse to a single "11:11" ------------------
$ws.Range("C12").Value = "11:11"
$ws.Range("D12").Value = "11:11"
$ws.Range("E12").Value = "11:11"
$ws.Range("F12").Value = "11:11"

# Clear the stray UNDERTIME hours value on row 12
$ws.Range("G12").ClearContents()

# --- Row 27 loses its day-number label ---------------------------------
$ws.Range("B27").ClearContents()

# --- Remove the extra day rows (old rows 28-49, days 17-37 + blank) ----
$ws.Rows("28:49").Delete()
